# Auto-generated Excel COM-interop script to update Leviathan Profits market data
# Applies the per-cell value changes derived from the target diff, sheet by sheet.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Cells.Item(74, 8).Value = 3675.7058  # H74: 3924.1875 -> 3675.7058
$ws.Cells.Item(74, 9).Value = 2771.5454  # I74: 3078.7 -> 2771.5454
$ws.Cells.Item(74, 11).Value = 2771.5454  # K74: 3078.7 -> 2771.5454
$ws.Cells.Item(74, 13).Value = -1835.5454  # M74: -2142.7 -> -1835.5454
# Row 77
$ws.Cells.Item(77, 8).Value = 3675.7058  # H77: 3924.1875 -> 3675.7058
$ws.Cells.Item(77, 9).Value = 2771.5454  # I77: 3078.7 -> 2771.5454
$ws.Cells.Item(77, 11).Value = 13857.727  # K77: 15393.5 -> 13857.727
$ws.Cells.Item(77, 13).Value = -9177.726999999999  # M77: -10713.5 -> -9177.726999999999
# Row 80
$ws.Cells.Item(80, 8).Value = 686.1667  # H80: 685.4 -> 686.1667
$ws.Cells.Item(80, 10).Value = 822  # J80: 866 -> 822
$ws.Cells.Item(80, 12).Value = 2466  # L80: 2598 -> 2466
$ws.Cells.Item(80, 14).Value = -4462  # N80: -4594 -> -4462
# Row 83
$ws.Cells.Item(83, 8).Value = 686.1667  # H83: 685.4 -> 686.1667
$ws.Cells.Item(83, 10).Value = 822  # J83: 866 -> 822
$ws.Cells.Item(83, 12).Value = 7398  # L83: 7794 -> 7398
$ws.Cells.Item(83, 14).Value = -17382  # N83: -17778 -> -17382
# Row 132
$ws.Cells.Item(132, 8).Value = 2906.5454  # H132: 2997.1904 -> 2906.5454
$ws.Cells.Item(132, 9).Value = 1558.8334  # I132: 1591.5294 -> 1558.8334
$ws.Cells.Item(132, 11).Value = 4676.5002  # K132: 4774.5882 -> 4676.5002
$ws.Cells.Item(132, 13).Value = -2146.5002  # M132: -2244.5882 -> -2146.5002
# Row 137
$ws.Cells.Item(137, 8).Value = 5037.933  # H137: 4766.8125 -> 5037.933
$ws.Cells.Item(137, 9).Value = 1776.6666  # I137: 1693.8462 -> 1776.6666
$ws.Cells.Item(137, 11).Value = 5329.9998  # K137: 5081.5386 -> 5329.9998
$ws.Cells.Item(137, 13).Value = -2779.9998  # M137: -2531.5386 -> -2779.9998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 28748.12  # H32: 28378.805 -> 28748.12
$ws.Cells.Item(32, 9).Value = 15634.479  # I32: 15415.571 -> 15634.479
$ws.Cells.Item(32, 11).Value = 15634.479  # K32: 15415.571 -> 15634.479
$ws.Cells.Item(32, 13).Value = -15347.479  # M32: -15128.571 -> -15347.479
# Row 74
$ws.Cells.Item(74, 8).Value = 1911.1428  # H74: 1670.6957 -> 1911.1428
$ws.Cells.Item(74, 9).Value = 1840.8182  # I74: 1416.8948 -> 1840.8182
$ws.Cells.Item(74, 10).Value = 2169  # J74: 2876.25 -> 2169
$ws.Cells.Item(74, 11).Value = 1840.8182  # K74: 1416.8948 -> 1840.8182
$ws.Cells.Item(74, 12).Value = 2169  # L74: 2876.25 -> 2169
$ws.Cells.Item(74, 13).Value = -966.8181999999999  # M74: -542.8948 -> -966.8181999999999
$ws.Cells.Item(74, 14).Value = -3917  # N74: -4624.25 -> -3917
# Row 77
$ws.Cells.Item(77, 8).Value = 1911.1428  # H77: 1670.6957 -> 1911.1428
$ws.Cells.Item(77, 9).Value = 1840.8182  # I77: 1416.8948 -> 1840.8182
$ws.Cells.Item(77, 10).Value = 2169  # J77: 2876.25 -> 2169
$ws.Cells.Item(77, 11).Value = 9204.091  # K77: 7084.474 -> 9204.091
$ws.Cells.Item(77, 12).Value = 10845  # L77: 14381.25 -> 10845
$ws.Cells.Item(77, 13).Value = -4836.091  # M77: -2716.474 -> -4836.091
$ws.Cells.Item(77, 14).Value = -19581  # N77: -23117.25 -> -19581
# Row 88
$ws.Cells.Item(88, 8).Value = 28572062  # H88: 50001176 -> 28572062
$ws.Cells.Item(88, 9).Value = 549.5  # I88: 1999 -> 549.5
$ws.Cells.Item(88, 10).Value = 40000668  # J88: 66667570 -> 40000668
$ws.Cells.Item(88, 11).Value = 549.5  # K88: 1999 -> 549.5
$ws.Cells.Item(88, 12).Value = 40000668  # L88: 66667570 -> 40000668
$ws.Cells.Item(88, 13).Value = -143.5  # M88: -1593 -> -143.5
$ws.Cells.Item(88, 14).Value = -40001480  # N88: -66668382 -> -40001480
# Row 91
$ws.Cells.Item(91, 8).Value = 28572062  # H91: 50001176 -> 28572062
$ws.Cells.Item(91, 9).Value = 549.5  # I91: 1999 -> 549.5
$ws.Cells.Item(91, 10).Value = 40000668  # J91: 66667570 -> 40000668
$ws.Cells.Item(91, 11).Value = 549.5  # K91: 1999 -> 549.5
$ws.Cells.Item(91, 12).Value = 40000668  # L91: 66667570 -> 40000668
$ws.Cells.Item(91, 13).Value = 854.5  # M91: -595 -> 854.5
$ws.Cells.Item(91, 14).Value = -40003476  # N91: -66670378 -> -40003476
# Row 102
$ws.Cells.Item(102, 8).Value = 102920.7  # H102: 69027.2 -> 102920.7
$ws.Cells.Item(102, 9).Value = 127401  # I102: 73600.57000000001 -> 127401
$ws.Cells.Item(102, 10).Value = 4999.5  # J102: 5000 -> 4999.5
$ws.Cells.Item(102, 11).Value = 127401  # K102: 73600.57000000001 -> 127401
$ws.Cells.Item(102, 12).Value = 4999.5  # L102: 5000 -> 4999.5
$ws.Cells.Item(102, 13).Value = -125779  # M102: -71978.57000000001 -> -125779
$ws.Cells.Item(102, 14).Value = -8243.5  # N102: -8244 -> -8243.5
# Row 122
$ws.Cells.Item(122, 8).Value = 2333.077  # H122: 2139.4119 -> 2333.077
$ws.Cells.Item(122, 10).Value = 2929.3333  # J122: 2118.2856 -> 2929.3333
$ws.Cells.Item(122, 12).Value = 8787.999899999999  # L122: 6354.8568 -> 8787.999899999999
$ws.Cells.Item(122, 14).Value = -13687.9999  # N122: -11254.8568 -> -13687.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Cells.Item(82, 8).Value = 16500.777  # H82: 15438.375 -> 16500.777
# Row 85
$ws.Cells.Item(85, 8).Value = 16500.777  # H85: 15438.375 -> 16500.777

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1646.8  # H31: 1711.7391 -> 1646.8
$ws.Cells.Item(31, 9).Value = 1269.2632  # I31: 1274.4736 -> 1269.2632
$ws.Cells.Item(31, 10).Value = 2842.3333  # J31: 3788.75 -> 2842.3333
$ws.Cells.Item(31, 11).Value = 1269.2632  # K31: 1274.4736 -> 1269.2632
$ws.Cells.Item(31, 12).Value = 2842.3333  # L31: 3788.75 -> 2842.3333
$ws.Cells.Item(31, 13).Value = -974.2632000000001  # M31: -979.4736 -> -974.2632000000001
$ws.Cells.Item(31, 14).Value = -3432.3333  # N31: -4378.75 -> -3432.3333
# Row 34
$ws.Cells.Item(34, 8).Value = 1646.8  # H34: 1711.7391 -> 1646.8
$ws.Cells.Item(34, 9).Value = 1269.2632  # I34: 1274.4736 -> 1269.2632
$ws.Cells.Item(34, 10).Value = 2842.3333  # J34: 3788.75 -> 2842.3333
$ws.Cells.Item(34, 11).Value = 1269.2632  # K34: 1274.4736 -> 1269.2632
$ws.Cells.Item(34, 12).Value = 2842.3333  # L34: 3788.75 -> 2842.3333
$ws.Cells.Item(34, 13).Value = -1067.2632  # M34: -1072.4736 -> -1067.2632
$ws.Cells.Item(34, 14).Value = -3246.3333  # N34: -4192.75 -> -3246.3333
# Row 58
$ws.Cells.Item(58, 8).Value = 1254.5555  # H58: 1228.9 -> 1254.5555
$ws.Cells.Item(58, 9).Value = 1254.5555  # I58: 1228.9 -> 1254.5555
$ws.Cells.Item(58, 11).Value = 1254.5555  # K58: 1228.9 -> 1254.5555
$ws.Cells.Item(58, 13).Value = -1051.5555  # M58: -1025.9 -> -1051.5555
# Row 74
$ws.Cells.Item(74, 8).Value = 20000  # H74: 30333.334 -> 20000
$ws.Cells.Item(74, 10).Value = 0  # J74: 35500 -> 0
$ws.Cells.Item(74, 12).Value = 0  # L74: 35500 -> 0
$ws.Cells.Item(74, 14).ClearContents()  # N74: was -37248
# Row 77
$ws.Cells.Item(77, 8).Value = 20000  # H77: 30333.334 -> 20000
$ws.Cells.Item(77, 10).Value = 0  # J77: 35500 -> 0
$ws.Cells.Item(77, 12).Value = 0  # L77: 106500 -> 0
$ws.Cells.Item(77, 14).ClearContents()  # N77: was -115236
# Row 86
$ws.Cells.Item(86, 8).Value = 9795  # H86: 9608 -> 9795
$ws.Cells.Item(86, 9).Value = 9158.333000000001  # I86: 8549.333000000001 -> 9158.333000000001
$ws.Cells.Item(86, 10).Value = 10750  # J86: 10666.667 -> 10750
$ws.Cells.Item(86, 11).Value = 9158.333000000001  # K86: 8549.333000000001 -> 9158.333000000001
$ws.Cells.Item(86, 12).Value = 10750  # L86: 10666.667 -> 10750
$ws.Cells.Item(86, 13).Value = -8035.333000000001  # M86: -7426.333000000001 -> -8035.333000000001
$ws.Cells.Item(86, 14).Value = -12996  # N86: -12912.667 -> -12996
# Row 89
$ws.Cells.Item(89, 8).Value = 9795  # H89: 9608 -> 9795
$ws.Cells.Item(89, 9).Value = 9158.333000000001  # I89: 8549.333000000001 -> 9158.333000000001
$ws.Cells.Item(89, 10).Value = 10750  # J89: 10666.667 -> 10750
$ws.Cells.Item(89, 11).Value = 45791.665  # K89: 42746.665 -> 45791.665
$ws.Cells.Item(89, 12).Value = 53750  # L89: 53333.335 -> 53750
$ws.Cells.Item(89, 13).Value = -40175.665  # M89: -37130.665 -> -40175.665
$ws.Cells.Item(89, 14).Value = -64982  # N89: -64565.335 -> -64982
# Row 122
$ws.Cells.Item(122, 8).Value = 4215.696  # H122: 4228.7393 -> 4215.696
$ws.Cells.Item(122, 9).Value = 3620.4443  # I122: 3637.111 -> 3620.4443
$ws.Cells.Item(122, 11).Value = 10861.3329  # K122: 10911.333 -> 10861.3329
$ws.Cells.Item(122, 13).Value = -8411.332900000001  # M122: -8461.332999999999 -> -8411.332900000001
# Row 136
$ws.Cells.Item(136, 8).Value = 1254.5555  # H136: 1228.9 -> 1254.5555
$ws.Cells.Item(136, 9).Value = 1254.5555  # I136: 1228.9 -> 1254.5555
$ws.Cells.Item(136, 11).Value = 3763.6665  # K136: 3686.7 -> 3763.6665
$ws.Cells.Item(136, 13).Value = -1213.6665  # M136: -1136.7 -> -1213.6665

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 623.75  # H5: 695.875 -> 623.75
$ws.Cells.Item(5, 9).Value = 623.75  # I5: 607 -> 623.75
$ws.Cells.Item(5, 10).Value = 0  # J5: 844 -> 0
$ws.Cells.Item(5, 11).Value = 1871.25  # K5: 1821 -> 1871.25
$ws.Cells.Item(5, 12).Value = 0  # L5: 2532 -> 0
$ws.Cells.Item(5, 13).Value = -1759.25  # M5: -1709 -> -1759.25
$ws.Cells.Item(5, 14).ClearContents()  # N5: was -2756
# Row 75
$ws.Cells.Item(75, 8).Value = 7292.385  # H75: 6105.8823 -> 7292.385
$ws.Cells.Item(75, 9).Value = 3478.75  # I75: 2416.4285 -> 3478.75
$ws.Cells.Item(75, 10).Value = 8987.333000000001  # J75: 8688.5 -> 8987.333000000001
$ws.Cells.Item(75, 11).Value = 10436.25  # K75: 7249.2855 -> 10436.25
$ws.Cells.Item(75, 12).Value = 26961.999  # L75: 26065.5 -> 26961.999
$ws.Cells.Item(75, 13).Value = -9438.25  # M75: -6251.2855 -> -9438.25
$ws.Cells.Item(75, 14).Value = -28957.999  # N75: -28061.5 -> -28957.999
# Row 78
$ws.Cells.Item(78, 8).Value = 7292.385  # H78: 6105.8823 -> 7292.385
$ws.Cells.Item(78, 9).Value = 3478.75  # I78: 2416.4285 -> 3478.75
$ws.Cells.Item(78, 10).Value = 8987.333000000001  # J78: 8688.5 -> 8987.333000000001
$ws.Cells.Item(78, 11).Value = 31308.75  # K78: 21747.8565 -> 31308.75
$ws.Cells.Item(78, 12).Value = 80885.997  # L78: 78196.5 -> 80885.997
$ws.Cells.Item(78, 13).Value = -26316.75  # M78: -16755.8565 -> -26316.75
$ws.Cells.Item(78, 14).Value = -90869.997  # N78: -88180.5 -> -90869.997
# Row 88
$ws.Cells.Item(88, 8).Value = 11666.667  # H88: 11249.75 -> 11666.667
$ws.Cells.Item(88, 10).Value = 11666.667  # J88: 11249.75 -> 11666.667
$ws.Cells.Item(88, 12).Value = 35000.001  # L88: 33749.25 -> 35000.001
$ws.Cells.Item(88, 14).Value = -35856.001  # N88: -34605.25 -> -35856.001
# Row 91
$ws.Cells.Item(91, 8).Value = 11666.667  # H91: 11249.75 -> 11666.667
$ws.Cells.Item(91, 10).Value = 11666.667  # J91: 11249.75 -> 11666.667
$ws.Cells.Item(91, 12).Value = 35000.001  # L91: 33749.25 -> 35000.001
$ws.Cells.Item(91, 14).Value = -37964.001  # N91: -36713.25 -> -37964.001
# Row 113
$ws.Cells.Item(113, 8).Value = 500.30768  # H113: 587.55554 -> 500.30768
$ws.Cells.Item(113, 10).Value = 454.9091  # J113: 541.1429000000001 -> 454.9091
$ws.Cells.Item(113, 12).Value = 1364.7273  # L113: 1623.4287 -> 1364.7273
$ws.Cells.Item(113, 14).Value = -5704.7273  # N113: -5963.4287 -> -5704.7273
# Row 128
$ws.Cells.Item(128, 8).Value = 399600.75  # H128: 399607.75 -> 399600.75
$ws.Cells.Item(128, 9).Value = 399600.75  # I128: 399607.75 -> 399600.75
$ws.Cells.Item(128, 11).Value = 1198802.25  # K128: 1198823.25 -> 1198802.25
$ws.Cells.Item(128, 13).Value = -1193822.25  # M128: -1193843.25 -> -1193822.25
# Row 134
$ws.Cells.Item(134, 8).Value = 6542.3335  # H134: 6831.4116 -> 6542.3335
$ws.Cells.Item(134, 9).Value = 2313.6667  # I134: 2376 -> 2313.6667
$ws.Cells.Item(134, 11).Value = 6941.000100000001  # K134: 7128 -> 6941.000100000001
$ws.Cells.Item(134, 13).Value = -1871.000100000001  # M134: -2058 -> -1871.000100000001
# Row 135
$ws.Cells.Item(135, 8).Value = 623.75  # H135: 695.875 -> 623.75
$ws.Cells.Item(135, 9).Value = 623.75  # I135: 607 -> 623.75
$ws.Cells.Item(135, 10).Value = 0  # J135: 844 -> 0
$ws.Cells.Item(135, 11).Value = 5613.75  # K135: 5463 -> 5613.75
$ws.Cells.Item(135, 12).Value = 0  # L135: 7596 -> 0
$ws.Cells.Item(135, 13).Value = -3078.75  # M135: -2928 -> -3078.75
$ws.Cells.Item(135, 14).ClearContents()  # N135: was -12666

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 103
$ws.Cells.Item(103, 8).Value = 40301  # H103: 0 -> 40301
$ws.Cells.Item(103, 10).Value = 40301  # J103: 0 -> 40301
$ws.Cells.Item(103, 12).Value = 40301  # L103: 0 -> 40301
$ws.Cells.Item(103, 14).Value = -42645  # N103: None -> -42645
# Row 122
$ws.Cells.Item(122, 8).Value = 144235.25  # H122: 136560.62 -> 144235.25
$ws.Cells.Item(122, 9).Value = 144235.25  # I122: 148454.4 -> 144235.25
$ws.Cells.Item(122, 10).Value = 0  # J122: 1764.3334 -> 0
$ws.Cells.Item(122, 11).Value = 432705.75  # K122: 445363.2 -> 432705.75
$ws.Cells.Item(122, 12).Value = 0  # L122: 5293.0002 -> 0
$ws.Cells.Item(122, 13).Value = -430255.75  # M122: -442913.2 -> -430255.75
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -10193.0002
# Row 126
$ws.Cells.Item(126, 8).Value = 3889.1667  # H126: 3762.1428 -> 3889.1667
$ws.Cells.Item(126, 9).Value = 3758.75  # I126: 4012 -> 3758.75
$ws.Cells.Item(126, 10).Value = 4150  # J126: 3574.75 -> 4150
$ws.Cells.Item(126, 11).Value = 11276.25  # K126: 12036 -> 11276.25
$ws.Cells.Item(126, 12).Value = 12450  # L126: 10724.25 -> 12450
$ws.Cells.Item(126, 13).Value = -8806.25  # M126: -9566 -> -8806.25
$ws.Cells.Item(126, 14).Value = -17390  # N126: -15664.25 -> -17390

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 2820.3845  # H68: 3055.5 -> 2820.3845
$ws.Cells.Item(68, 9).Value = 2575.111  # I68: 2811.1428 -> 2575.111
$ws.Cells.Item(68, 10).Value = 3372.25  # J68: 3397.6 -> 3372.25
$ws.Cells.Item(68, 11).Value = 2575.111  # K68: 2811.1428 -> 2575.111
$ws.Cells.Item(68, 12).Value = 3372.25  # L68: 3397.6 -> 3372.25
$ws.Cells.Item(68, 13).Value = -1826.111  # M68: -2062.1428 -> -1826.111
$ws.Cells.Item(68, 14).Value = -4870.25  # N68: -4895.6 -> -4870.25
# Row 71
$ws.Cells.Item(71, 8).Value = 2820.3845  # H71: 3055.5 -> 2820.3845
$ws.Cells.Item(71, 9).Value = 2575.111  # I71: 2811.1428 -> 2575.111
$ws.Cells.Item(71, 10).Value = 3372.25  # J71: 3397.6 -> 3372.25
$ws.Cells.Item(71, 11).Value = 12875.555  # K71: 14055.714 -> 12875.555
$ws.Cells.Item(71, 12).Value = 16861.25  # L71: 16988 -> 16861.25
$ws.Cells.Item(71, 13).Value = -9131.555  # M71: -10311.714 -> -9131.555
$ws.Cells.Item(71, 14).Value = -24349.25  # N71: -24476 -> -24349.25
# Row 82
$ws.Cells.Item(82, 8).Value = 2029.7778  # H82: 2259.8 -> 2029.7778
$ws.Cells.Item(82, 9).Value = 1657.8  # I82: 1150 -> 1657.8
$ws.Cells.Item(82, 10).Value = 2494.75  # J82: 2999.6667 -> 2494.75
$ws.Cells.Item(82, 11).Value = 1657.8  # K82: 1150 -> 1657.8
$ws.Cells.Item(82, 12).Value = 2494.75  # L82: 2999.6667 -> 2494.75
$ws.Cells.Item(82, 13).Value = -1296.8  # M82: -789 -> -1296.8
$ws.Cells.Item(82, 14).Value = -3216.75  # N82: -3721.6667 -> -3216.75
# Row 85
$ws.Cells.Item(85, 8).Value = 2029.7778  # H85: 2259.8 -> 2029.7778
$ws.Cells.Item(85, 9).Value = 1657.8  # I85: 1150 -> 1657.8
$ws.Cells.Item(85, 10).Value = 2494.75  # J85: 2999.6667 -> 2494.75
$ws.Cells.Item(85, 11).Value = 1657.8  # K85: 1150 -> 1657.8
$ws.Cells.Item(85, 12).Value = 2494.75  # L85: 2999.6667 -> 2494.75
$ws.Cells.Item(85, 13).Value = -409.8  # M85: 98 -> -409.8
$ws.Cells.Item(85, 14).Value = -4990.75  # N85: -5495.6667 -> -4990.75
# Row 132
$ws.Cells.Item(132, 8).Value = 717716.5  # H132: 520014.38 -> 717716.5
$ws.Cells.Item(132, 9).Value = 1073575.2  # I132: 683564.75 -> 1073575.2
$ws.Cells.Item(132, 11).Value = 3220725.6  # K132: 2050694.25 -> 3220725.6
$ws.Cells.Item(132, 13).Value = -3218195.6  # M132: -2048164.25 -> -3218195.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Cells.Item(136, 8).Value = 3110.543  # H136: 3270.2424 -> 3110.543
$ws.Cells.Item(136, 9).Value = 3229.4  # I136: 3460.4075 -> 3229.4
$ws.Cells.Item(136, 10).Value = 2397.4  # J136: 2414.5 -> 2397.4
$ws.Cells.Item(136, 11).Value = 9688.200000000001  # K136: 10381.2225 -> 9688.200000000001
$ws.Cells.Item(136, 12).Value = 7192.200000000001  # L136: 7243.5 -> 7192.200000000001
$ws.Cells.Item(136, 13).Value = -7138.200000000001  # M136: -7831.2225 -> -7138.200000000001
$ws.Cells.Item(136, 14).Value = -12292.2  # N136: -12343.5 -> -12292.2
